$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-21, columns A-C (row 1 header is unchanged)
$data = @{
    2  = @("平潭发展", "三六零",   "平潭发展")
    3  = @("福龙马",   "平潭发展", "福龙马")
    4  = @("神州信息", "东方精工", "东方精工")
    5  = @("东方精工", "福龙马",   "神州信息")
    6  = @("三六零",   "特变电工", "特变电工")
    7  = @("海马汽车", "神州信息", "山子高科")
    8  = @("特变电工", "工业富联", "海南发展")
    9  = @("粤 传 媒", "山子高科", "华胜天成")
    10 = @("山子高科", "包钢股份", "三六零")
    11 = @("兰石重装", "粤 传 媒", "海马汽车")
    12 = @("工业富联", "吉视传媒", "凯美特气")
    13 = @("华胜天成", "上海电气", "粤传媒")
    14 = @("万向钱潮", "东方明珠", "三花智控")
    15 = @("吉视传媒", "兰石重装", "太极实业")
    16 = @("包钢股份", "中核科技", "阳光电源")
    17 = @("太极实业", "华胜天成", "吉视传媒")
    18 = @("三花智控", "海马汽车", "工业富联")
    19 = @("东方明珠", "太极实业", "合富中国")
    20 = @("阳光电源", "航天智装", "万向钱潮")
    21 = @("美瑞新材", "天际股份", "亚太药业")
}

foreach ($rowNum in $data.Keys) {
    $vals = $data[$rowNum]
    $ws.Cells.Item($rowNum, 1).Value = $vals[0]
    $ws.Cells.Item($rowNum, 2).Value = $vals[1]
    $ws.Cells.Item($rowNum, 3).Value = $vals[2]
}
